$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    4329.98560929835,
    4208.331502926498,
    4036.612779150249,
    3938.223653360273,
    3938.223653360273,
    3938.223653360273,
    3938.223653360273,
    3872.912727445029,
    3872.912727445029,
    3872.912727445029,
    3872.912727445029
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
